# The diff shows a new data row being inserted as row 73 ("weekly" price
# observation added at the top of the "Cebollín baby" series), pushing all
# the existing data rows 73:107 down by one (to 74:108). The dimension grows
# from A1:R107 to A1:R108.
#
# To reproduce this with Excel COM semantics:
#   1. Duplicate row 73 (Copy + Insert) so that the new row 73 starts out as
#      an exact clone of the original row 73 - this pushes the original
#      row 73 (and everything below it) down to row 74, preserving all
#      values/number formats/styles automatically.
#   2. Overwrite the four cells that differ for the new observation
#      (Fecha/D, Precio mínimo/K, Precio promedio ponderado/M, Precio $/Kg/P).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(73).Copy()
$ws.Rows.Item(73).Insert()

$ws.Range("D73").Value = 44845
$ws.Range("K73").Value = 1300
$ws.Range("M73").Value = 1400
$ws.Range("P73").Value = 700
